$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (date bumped a day, dose number changed)
$ws.Range("A2").Value = "14/08/2023"
$ws.Range("B2").Value = "516.589.644"

# Remove row 3 entirely (was A3:B3) -- shrinks used range to A1:B2
$ws.Rows.Item(3).Delete()

# Set explicit column widths for A and B so the saved <col> width reads "13"
# (Excel's COM ColumnWidth is offset from the stored OOXML width by 5/6,
# i.e. stored = ColumnWidth + 5/6, so 13 - 5/6 = 12.1666... gives width="13")
$ws.Range("A:B").ColumnWidth = 12.166666666666666
